# Update countries & provincias Spain
# Applies the data refresh described in the commit:
#  - Moves "Republica de Yibuti" earlier in the shared-strings order
#    (right after "Guayana Francesa"), which shifts the country labels
#    for the Jamaica/Barbados/Yibuti block of rows by one position and
#    refreshes their COVID figures.
#  - Updates several countries' case/death counters (Alemania, Austria,
#    Pakistan, Libano, Banglades).
#  - Refreshes the "Datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 15:22"

# --- Row 7: Alemania ---
$ws.Range("B7").Value = 97052
$ws.Range("C7").Value = 960
$ws.Range("E7").Value = 69174
$ws.Range("G7").Value = 34
$ws.Range("H7").Value = 1478

# --- Row 17: Austria ---
$ws.Range("B17").Value = 11907
$ws.Range("C17").Value = 126
$ws.Range("E17").Value = 8705

# --- Row 37: Pakistan ---
$ws.Range("B37").Value = 2899
$ws.Range("C37").Value = 81
$ws.Range("E37").Value = 2646

# --- Row 81: Libano ---
$ws.Range("F81").Value = 28

# --- Row 128: Banglades ---
$ws.Range("D128").Value = 33
$ws.Range("E128").Value = 46

# --- Rows 136-138: Republica de Yibuti inserted before Jamaica/Barbados ---
# Row 136 becomes Republica de Yibuti with refreshed figures
$ws.Range("A136").Value = "Republica de Yibuti"
$ws.Range("B136").Value = 59
$ws.Range("C136").Value = 9
$ws.Range("D136").Value = 9
$ws.Range("E136").Value = 50
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 0

# Row 137 becomes Jamaica (previous Jamaica figures)
$ws.Range("A137").Value = "Jamaica"
$ws.Range("B137").Value = 55
$ws.Range("C137").Value = 2
$ws.Range("D137").Value = 7
$ws.Range("E137").Value = 45
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 3

# Row 138 becomes Barbados (previous Barbados figures)
$ws.Range("A138").Value = "Barbados"
$ws.Range("B138").Value = 52
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 0
$ws.Range("E138").Value = 52
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 0
